$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Delivery_results")

# New headers for columns F and G, matching the header style used by A1:E1
$ws.Range("F1").Value = "Average_waiting_time_(minutes)"
$ws.Range("G1").Value = "Average_queue_length"

$ws.Range("E1").Copy()
$ws.Range("F1:G1").PasteSpecial(-4122)

# Update existing column D and E values, and add F/G values per row
$ws.Range("D2").Value = 55
$ws.Range("E2").Value = 88.70999999999999
$ws.Range("F2").Value = 1.8
$ws.Range("G2").Value = 6.67

$ws.Range("D3").Value = 32
$ws.Range("E3").Value = 53.33
$ws.Range("F3").Value = 0.28
$ws.Range("G3").Value = 1.53

$ws.Range("D4").Value = 43
$ws.Range("E4").Value = 100
$ws.Range("F4").Value = 0.65
$ws.Range("G4").Value = 3.77

$ws.Range("D5").Value = 15
$ws.Range("E5").Value = 83.33
$ws.Range("F5").Value = 0.88
$ws.Range("G5").Value = 4.8

$ws.Range("D6").Value = 37
$ws.Range("E6").Value = 100
$ws.Range("F6").Value = 0.71
$ws.Range("G6").Value = 3.43
